# Applies the updated coin price / volume figures (and the HuobiToken /
# ImmutableX row swap) described by the commit's XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new text, and whether the new text must be forced
# to stay a STRING (many "Price" values look like numbers, e.g. 0.1270 or
# 0.000009354 - plain .Value assignment would let Excel reinterpret them as
# numeric and silently drop the trailing zero / switch to scientific notation).
$updates = @(
    @{ Cell = "D2"; Value = '28.946.43'; AsText = $false },
    @{ Cell = "E2"; Value = '  -0.92%  '; AsText = $false },
    @{ Cell = "D3"; Value = '1.814.29'; AsText = $false },
    @{ Cell = "E3"; Value = '  -0.90%  '; AsText = $false },
    @{ Cell = "D4"; Value = '1.002'; AsText = $true },
    @{ Cell = "E4"; Value = '  +0.18%  '; AsText = $false },
    @{ Cell = "D5"; Value = '232.45'; AsText = $true },
    @{ Cell = "E5"; Value = '  -2.11%  '; AsText = $false },
    @{ Cell = "D6"; Value = '0.5911'; AsText = $true },
    @{ Cell = "E6"; Value = '  -3.03%  '; AsText = $false },
    @{ Cell = "E7"; Value = '  +0.24%  '; AsText = $false },
    @{ Cell = "D8"; Value = '0.2747'; AsText = $true },
    @{ Cell = "E8"; Value = '  -2.52%  '; AsText = $false },
    @{ Cell = "D9"; Value = '0.06756'; AsText = $true },
    @{ Cell = "E9"; Value = '  -4.72%  '; AsText = $false },
    @{ Cell = "D10"; Value = '22.91'; AsText = $true },
    @{ Cell = "E10"; Value = '  -4.10%  '; AsText = $false },
    @{ Cell = "D11"; Value = '0.07496'; AsText = $true },
    @{ Cell = "E11"; Value = '  -1.99%  '; AsText = $false },
    @{ Cell = "D12"; Value = '1.871.13'; AsText = $false },
    @{ Cell = "E12"; Value = '  +2.04%  '; AsText = $false },
    @{ Cell = "D13"; Value = '4.674'; AsText = $true },
    @{ Cell = "E13"; Value = '  -2.91%  '; AsText = $false },
    @{ Cell = "D14"; Value = '0.6237'; AsText = $true },
    @{ Cell = "E14"; Value = '  -1.50%  '; AsText = $false },
    @{ Cell = "D15"; Value = '0.000009354'; AsText = $true },
    @{ Cell = "E15"; Value = '  -6.23%  '; AsText = $false },
    @{ Cell = "D16"; Value = '74.67'; AsText = $true },
    @{ Cell = "E16"; Value = '  -6.10%  '; AsText = $false },
    @{ Cell = "D17"; Value = '28.719.86'; AsText = $false },
    @{ Cell = "E17"; Value = '  -1.74%  '; AsText = $false },
    @{ Cell = "D18"; Value = '5.432'; AsText = $true },
    @{ Cell = "E18"; Value = '  -8.97%  '; AsText = $false },
    @{ Cell = "E19"; Value = '  +0.25%  '; AsText = $false },
    @{ Cell = "D20"; Value = '208.25'; AsText = $true },
    @{ Cell = "E20"; Value = '  -9.26%  '; AsText = $false },
    @{ Cell = "D21"; Value = '11.38'; AsText = $true },
    @{ Cell = "E21"; Value = '  -3.76%  '; AsText = $false },
    @{ Cell = "D22"; Value = '6.772'; AsText = $true },
    @{ Cell = "E22"; Value = '  -3.89%  '; AsText = $false },
    @{ Cell = "D23"; Value = '1.004'; AsText = $true },
    @{ Cell = "E23"; Value = '  +0.11%  '; AsText = $false },
    @{ Cell = "E24"; Value = '  -0.63%  '; AsText = $false },
    @{ Cell = "D25"; Value = '0.1270'; AsText = $true },
    @{ Cell = "E25"; Value = '  -2.44%  '; AsText = $false },
    @{ Cell = "D26"; Value = '7.790'; AsText = $true },
    @{ Cell = "E26"; Value = '  -4.05%  '; AsText = $false },
    @{ Cell = "D27"; Value = '16.30'; AsText = $true },
    @{ Cell = "E27"; Value = '  -2.61%  '; AsText = $false },
    @{ Cell = "D28"; Value = '0.06361'; AsText = $true },
    @{ Cell = "E28"; Value = '  -5.59%  '; AsText = $false },
    @{ Cell = "D29"; Value = '1.401'; AsText = $true },
    @{ Cell = "E29"; Value = '  -5.52%  '; AsText = $false },
    @{ Cell = "D30"; Value = '1.429'; AsText = $true },
    @{ Cell = "E30"; Value = '  -2.13%  '; AsText = $false },
    @{ Cell = "D31"; Value = '3.728'; AsText = $true },
    @{ Cell = "E31"; Value = '  -3.03%  '; AsText = $false },
    @{ Cell = "D32"; Value = '3.680'; AsText = $true },
    @{ Cell = "E32"; Value = '  -4.32%  '; AsText = $false },
    @{ Cell = "D33"; Value = '1.687'; AsText = $true },
    @{ Cell = "E33"; Value = '  -2.85%  '; AsText = $false },
    @{ Cell = "D34"; Value = '1.048'; AsText = $true },
    @{ Cell = "E34"; Value = '  -7.48%  '; AsText = $false },
    @{ Cell = "B35"; Value = 'ImmutableX'; AsText = $false },
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; AsText = $false },
    @{ Cell = "D35"; Value = '0.6336'; AsText = $true },
    @{ Cell = "E35"; Value = '  -3.32%  '; AsText = $false },
    @{ Cell = "B36"; Value = 'HuobiToken'; AsText = $false },
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; AsText = $false },
    @{ Cell = "D36"; Value = '2.521'; AsText = $true },
    @{ Cell = "E36"; Value = '  -1.24%  '; AsText = $false },
    @{ Cell = "D37"; Value = '2.729'; AsText = $true },
    @{ Cell = "E37"; Value = '  -1.19%  '; AsText = $false },
    @{ Cell = "D38"; Value = '6.422'; AsText = $true },
    @{ Cell = "E38"; Value = '  -2.63%  '; AsText = $false },
    @{ Cell = "E39"; Value = '  -4.35%  '; AsText = $false },
    @{ Cell = "D40"; Value = '1.133.00'; AsText = $false },
    @{ Cell = "E40"; Value = '  -8.39%  '; AsText = $false },
    @{ Cell = "D41"; Value = '0.8676'; AsText = $true },
    @{ Cell = "E41"; Value = '  -5.94%  '; AsText = $false },
    @{ Cell = "D42"; Value = '1.003'; AsText = $true },
    @{ Cell = "E42"; Value = '  +0.21%  '; AsText = $false },
    @{ Cell = "E43"; Value = '  -0.94%  '; AsText = $false },
    @{ Cell = "D44"; Value = '99.82'; AsText = $true },
    @{ Cell = "E44"; Value = '  -1.09%  '; AsText = $false },
    @{ Cell = "D45"; Value = '60.33'; AsText = $true },
    @{ Cell = "E45"; Value = '  -5.08%  '; AsText = $false },
    @{ Cell = "D46"; Value = '0.00000000113'; AsText = $true },
    @{ Cell = "E46"; Value = '  -2.48%  '; AsText = $false },
    @{ Cell = "D47"; Value = '1.573'; AsText = $true },
    @{ Cell = "E47"; Value = '  -3.54%  '; AsText = $false },
    @{ Cell = "D48"; Value = '0.05461'; AsText = $true },
    @{ Cell = "E48"; Value = '  -1.75%  '; AsText = $false },
    @{ Cell = "D49"; Value = '0.4509'; AsText = $true },
    @{ Cell = "E49"; Value = '  -1.19%  '; AsText = $false },
    @{ Cell = "D50"; Value = '8.247'; AsText = $true },
    @{ Cell = "E50"; Value = '  -4.08%  '; AsText = $false },
    @{ Cell = "E51"; Value = '  -0.22%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        # Force text storage, write the value, then drop back to the default
        # style so no stray number-format style is left on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}

